# Ajout de la formule pour calculer la mise optimale dans wtfSpredSheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New label / value cells (rows 11, 13-15)
$ws.Range("I11").Value = "Formule:"
$ws.Range("J11").Value = "(3^round - 1)/2"

$ws.Range("I13").Value = "Round"
$ws.Range("J13").Value = 8

$ws.Range("I14").Value = "Cash en banque"
$ws.Range("J14").Value = 44.2

$ws.Range("I15").Value = "Mise optimale"
$ws.Range("J15").Formula = "=_xlfn.FLOOR.MATH(J14/(((3^J13)-1)/2),0.0001)"
$ws.Range("J15").NumberFormat = "0.0000"

# E3 now derives the mise from the optimal-bet computation
$ws.Range("E3").Formula = "=J15"
$ws.Range("E3").NumberFormat = "0.0000"

# Move the active selection like the author's last save
$ws.Range("E4").Select()
